# Login tests and ProductsTests preparations
#
# - Fix the typo in the 4th sheet's name: "WrongUserName" -> "WrongUsername"
# - Move the active/selected tab from "Login" (sheet 1) to "WrongUsername"
#   (sheet 4), which also flips each sheet's tabSelected flag accordingly.

$wb = $excel.ActiveWorkbook

$wrongUserNameSheet = $wb.Worksheets.Item(4)
$wrongUserNameSheet.Name = "WrongUsername"

# Activating this sheet clears tabSelected on whichever sheet previously had
# it (Login) and sets tabSelected="1" / activeTab on this one instead.
$wrongUserNameSheet.Activate()
